$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp shown above the table
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Abril de 2020 a las 19:22"

# Refresh country case-count figures; the table is kept sorted by total
# cases (column B) descending, so several rows also change which country
# they display as case counts moved countries up/down in the ranking.
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 410916
$ws.Cells.Item(4, 3).Value = 10581
$ws.Cells.Item(4, 4).Value = 22081
$ws.Cells.Item(4, 5).Value = 374625
$ws.Cells.Item(4, 6).Value = 9220
$ws.Cells.Item(4, 7).Value = 1369
$ws.Cells.Item(4, 8).Value = 14210

$ws.Cells.Item(7, 1).Value = "Alemania"
$ws.Cells.Item(7, 2).Value = 109702
$ws.Cells.Item(7, 3).Value = 2039
$ws.Cells.Item(7, 4).Value = 36081
$ws.Cells.Item(7, 5).Value = 71516
$ws.Cells.Item(7, 6).Value = 4895
$ws.Cells.Item(7, 7).Value = 89
$ws.Cells.Item(7, 8).Value = 2105

$ws.Cells.Item(12, 1).Value = "Turquia"
$ws.Cells.Item(12, 2).Value = 38226
$ws.Cells.Item(12, 3).Value = 4117
$ws.Cells.Item(12, 4).Value = 1846
$ws.Cells.Item(12, 5).Value = 35568
$ws.Cells.Item(12, 6).Value = 1492
$ws.Cells.Item(12, 7).Value = 87
$ws.Cells.Item(12, 8).Value = 812

$ws.Cells.Item(14, 1).Value = "Suiza"
$ws.Cells.Item(14, 2).Value = 23248
$ws.Cells.Item(14, 3).Value = 995
$ws.Cells.Item(14, 4).Value = 9800
$ws.Cells.Item(14, 5).Value = 12553
$ws.Cells.Item(14, 6).Value = 391
$ws.Cells.Item(14, 7).Value = 74
$ws.Cells.Item(14, 8).Value = 895

$ws.Cells.Item(17, 1).Value = "Brasil"
$ws.Cells.Item(17, 2).Value = 14324
$ws.Cells.Item(17, 3).Value = 290
$ws.Cells.Item(17, 4).Value = 127
$ws.Cells.Item(17, 5).Value = 13491
$ws.Cells.Item(17, 6).Value = 296
$ws.Cells.Item(17, 7).Value = 20
$ws.Cells.Item(17, 8).Value = 706

$ws.Cells.Item(19, 1).Value = "Austria"
$ws.Cells.Item(19, 2).Value = 12920
$ws.Cells.Item(19, 3).Value = 281
$ws.Cells.Item(19, 4).Value = 4512
$ws.Cells.Item(19, 5).Value = 8135
$ws.Cells.Item(19, 6).Value = 267
$ws.Cells.Item(19, 7).Value = 30
$ws.Cells.Item(19, 8).Value = 273

$ws.Cells.Item(24, 1).Value = "Noruega"
$ws.Cells.Item(24, 2).Value = 6086
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 32
$ws.Cells.Item(24, 5).Value = 5953
$ws.Cells.Item(24, 6).Value = 78
$ws.Cells.Item(24, 7).Value = 12
$ws.Cells.Item(24, 8).Value = 101

$ws.Cells.Item(25, 1).Value = "Irlanda"
$ws.Cells.Item(25, 2).Value = 6074
$ws.Cells.Item(25, 3).Value = 365
$ws.Cells.Item(25, 4).Value = 25
$ws.Cells.Item(25, 5).Value = 5814
$ws.Cells.Item(25, 6).Value = 165
$ws.Cells.Item(25, 7).Value = 25
$ws.Cells.Item(25, 8).Value = 235

$ws.Cells.Item(26, 1).Value = "Australia"
$ws.Cells.Item(26, 2).Value = 6013
$ws.Cells.Item(26, 3).Value = 25
$ws.Cells.Item(26, 4).Value = 2813
$ws.Cells.Item(26, 5).Value = 3150
$ws.Cells.Item(26, 6).Value = 87
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 50

$ws.Cells.Item(27, 1).Value = "India"
$ws.Cells.Item(27, 2).Value = 5749
$ws.Cells.Item(27, 3).Value = 398
$ws.Cells.Item(27, 4).Value = 506
$ws.Cells.Item(27, 5).Value = 5065
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 18
$ws.Cells.Item(27, 8).Value = 178

$ws.Cells.Item(89, 1).Value = "Afganistan"
$ws.Cells.Item(89, 2).Value = 444
$ws.Cells.Item(89, 3).Value = 21
$ws.Cells.Item(89, 4).Value = 29
$ws.Cells.Item(89, 5).Value = 401
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 14

$ws.Cells.Item(90, 1).Value = "Uruguay"
$ws.Cells.Item(90, 2).Value = 424
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 150
$ws.Cells.Item(90, 5).Value = 267
$ws.Cells.Item(90, 6).Value = 14
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 7

$ws.Cells.Item(96, 1).Value = "Jordania"
$ws.Cells.Item(96, 2).Value = 358
$ws.Cells.Item(96, 3).Value = 5
$ws.Cells.Item(96, 4).Value = 150
$ws.Cells.Item(96, 5).Value = 202
$ws.Cells.Item(96, 6).Value = 5
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 6

$ws.Cells.Item(157, 1).Value = "Liberia"
$ws.Cells.Item(157, 2).Value = 31
$ws.Cells.Item(157, 3).Value = 17
$ws.Cells.Item(157, 4).Value = 3
$ws.Cells.Item(157, 5).Value = 24
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 1
$ws.Cells.Item(157, 8).Value = 4

$ws.Cells.Item(158, 1).Value = "Gabon"
$ws.Cells.Item(158, 2).Value = 30
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 1
$ws.Cells.Item(158, 5).Value = 28
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 1

$ws.Cells.Item(159, 1).Value = "Haiti"
$ws.Cells.Item(159, 2).Value = 27
$ws.Cells.Item(159, 3).Value = 2
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 26
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 1

$ws.Cells.Item(160, 1).Value = "Benin"
$ws.Cells.Item(160, 2).Value = 26
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 5
$ws.Cells.Item(160, 5).Value = 20
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 1

$ws.Cells.Item(161, 1).Value = "Tanzania"
$ws.Cells.Item(161, 2).Value = 25
$ws.Cells.Item(161, 3).Value = 1
$ws.Cells.Item(161, 4).Value = 5
$ws.Cells.Item(161, 5).Value = 19
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 1

$ws.Cells.Item(162, 1).Value = "Birmania"
$ws.Cells.Item(162, 2).Value = 22
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 0
$ws.Cells.Item(162, 5).Value = 19
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 2
$ws.Cells.Item(162, 8).Value = 3

$ws.Cells.Item(163, 1).Value = "Libia"
$ws.Cells.Item(163, 2).Value = 21
$ws.Cells.Item(163, 3).Value = 1
$ws.Cells.Item(163, 4).Value = 2
$ws.Cells.Item(163, 5).Value = 18
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 1

$ws.Cells.Item(164, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(164, 2).Value = 19
$ws.Cells.Item(164, 3).Value = 0
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 17
$ws.Cells.Item(164, 6).Value = 1
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 2

$ws.Cells.Item(165, 1).Value = "Siria"
$ws.Cells.Item(165, 2).Value = 19
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 3
$ws.Cells.Item(165, 5).Value = 14
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 2

$ws.Cells.Item(166, 1).Value = "Maldivas"
$ws.Cells.Item(166, 2).Value = 19
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 4).Value = 13
$ws.Cells.Item(166, 5).Value = 6
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 0

$ws.Cells.Item(167, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(167, 2).Value = 18
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 1
$ws.Cells.Item(167, 5).Value = 17
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 0
$ws.Cells.Item(167, 8).Value = 0

$ws.Cells.Item(168, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(168, 2).Value = 17
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = 17
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 0

$ws.Cells.Item(169, 1).Value = "Mozambique"
$ws.Cells.Item(169, 2).Value = 17
$ws.Cells.Item(169, 3).Value = 7
$ws.Cells.Item(169, 4).Value = 1
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0

$ws.Cells.Item(170, 1).Value = "Angola"
$ws.Cells.Item(170, 2).Value = 17
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 2
$ws.Cells.Item(170, 5).Value = 13
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 2

$ws.Cells.Item(172, 1).Value = "Namibia"
$ws.Cells.Item(172, 2).Value = 16
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(172, 4).Value = 3
$ws.Cells.Item(172, 5).Value = 13
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(173, 1).Value = "Mongolia"
$ws.Cells.Item(173, 2).Value = 16
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(173, 4).Value = 4
$ws.Cells.Item(173, 5).Value = 12
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(175, 1).Value = "Fiyi"
$ws.Cells.Item(175, 2).Value = 15
$ws.Cells.Item(175, 3).Value = 0
$ws.Cells.Item(175, 4).Value = 0
$ws.Cells.Item(175, 5).Value = 15
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

$ws.Cells.Item(176, 1).Value = "Dominica"
$ws.Cells.Item(176, 2).Value = 15
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 1
$ws.Cells.Item(176, 5).Value = 14
$ws.Cells.Item(176, 6).Value = 0
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = "Santa Lucia"
$ws.Cells.Item(177, 2).Value = 14
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 1
$ws.Cells.Item(177, 5).Value = 13
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 0
$ws.Cells.Item(177, 8).Value = 0

$ws.Cells.Item(178, 1).Value = "Sudan"
$ws.Cells.Item(178, 2).Value = 14
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 2
$ws.Cells.Item(178, 5).Value = 10
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 0
$ws.Cells.Item(178, 8).Value = 2

$ws.Cells.Item(183, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(183, 2).Value = 11
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 11
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

$ws.Cells.Item(184, 1).Value = "Seychelles"
$ws.Cells.Item(184, 2).Value = 11
$ws.Cells.Item(184, 3).Value = 0
$ws.Cells.Item(184, 4).Value = 0
$ws.Cells.Item(184, 5).Value = 11
$ws.Cells.Item(184, 6).Value = 0
$ws.Cells.Item(184, 7).Value = 0
$ws.Cells.Item(184, 8).Value = 0

$ws.Cells.Item(192, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(192, 2).Value = 8
$ws.Cells.Item(192, 3).Value = 0
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 5).Value = 7
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 1

$ws.Cells.Item(193, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(193, 2).Value = 8
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 5).Value = 7
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = "Malaui"
$ws.Cells.Item(194, 2).Value = 8
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 7
$ws.Cells.Item(194, 6).Value = 1
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 1

$ws.Cells.Item(195, 1).Value = "Belice"
$ws.Cells.Item(195, 2).Value = 8
$ws.Cells.Item(195, 3).Value = 1
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 7
$ws.Cells.Item(195, 6).Value = 1
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 1

$ws.Cells.Item(208, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(208, 2).Value = 3
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 0
$ws.Cells.Item(208, 5).Value = 3
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Anguila"
$ws.Cells.Item(209, 2).Value = 3
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 0
$ws.Cells.Item(209, 5).Value = 3
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(211, 1).Value = "Sudan del Sur"
$ws.Cells.Item(211, 2).Value = 2
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 0
$ws.Cells.Item(211, 5).Value = 2
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 2).Value = 2
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 0
$ws.Cells.Item(213, 5).Value = 2
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

